$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Saldo Anterior" section ---
$ws.Range("B2").Value = "R$100.00"
$ws.Range("A3").Value = "04-2024"

# --- Update Receitas (income) values ---
$ws.Range("B5").Value = "R$233.00"
$ws.Range("B6").Value = "R$233.00"
$ws.Range("B7").Value = "R$233.00"
$ws.Range("B8").Value = "R$233.00"
$ws.Range("B9").Value = "R$233.00"
$ws.Range("B10").Value = "R$233.00"
$ws.Range("B12").Value = "R$1398.00"

# --- Update Despesas (expense) labels/values before removing rows ---
$ws.Range("A15").Value = "copasa"
$ws.Range("B15").Value = "R$100.00"
$ws.Range("A16").Value = "skdks"
$ws.Range("B16").Value = "R$232.00"

# --- Delete rows 17 and 18 (the "Cera" and "Randup" expense line items) ---
# This shifts everything below up by two rows, matching the diff which
# removes those two rows entirely and shifts the rest of the sheet up.
$ws.Range("A17:D18").EntireRow.Delete()

# --- After the shift, set the new totals / summary values ---
# Row 18 is now "Total de Despesas:" (was row 20). It should keep the
# 23.25pt custom row height that the deleted "Randup" row (old row 18) had.
$ws.Range("B18").Value = "R$332.00"
$ws.Rows.Item(18).RowHeight = 23.25

# Row 21 is now "Saldo Mês anterior:" (was row 23)
$ws.Range("B21").Value = "R$100.00"

# Row 22 is now "Receitas Realizadas:" (was row 24)
$ws.Range("B22").Value = "R$1398.00"

# Row 23 is now "Despesas Realizadas" (was row 25)
$ws.Range("B23").Value = "R$332.00"

# Row 25 is now "Saldo Atual:" (was row 27)
$ws.Range("B25").Value = "R$1166.0"

$wb.Save()
